# Refresh the cryptocurrency price/volume table with newly scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some refreshed "Price" values look like plain numbers (e.g. "278.41").
# The source data stores these as literal text (matching the original
# worksheet, which keeps every Price cell as a text string), so each such
# cell is temporarily switched to Text format while its value is written,
# then restored to its original style to avoid altering formatting.
$origStyles = @{}
$origStyles["D4"] = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$origStyles["D5"] = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$origStyles["D6"] = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$origStyles["D7"] = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$origStyles["D8"] = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$origStyles["D9"] = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$origStyles["D11"] = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$origStyles["D12"] = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$origStyles["D13"] = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$origStyles["D14"] = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$origStyles["D15"] = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$origStyles["D16"] = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$origStyles["D17"] = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$origStyles["D19"] = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$origStyles["D20"] = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$origStyles["D21"] = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$origStyles["D22"] = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$origStyles["D23"] = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$origStyles["D25"] = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$origStyles["D26"] = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$origStyles["D27"] = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$origStyles["D28"] = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$origStyles["D30"] = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$origStyles["D31"] = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$origStyles["D32"] = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$origStyles["D33"] = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$origStyles["D34"] = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$origStyles["D36"] = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$origStyles["D37"] = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$origStyles["D38"] = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$origStyles["D39"] = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$origStyles["D40"] = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$origStyles["D42"] = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$origStyles["D43"] = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$origStyles["D44"] = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$origStyles["D45"] = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$origStyles["D46"] = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$origStyles["D47"] = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$origStyles["D48"] = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$origStyles["D49"] = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$origStyles["D50"] = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$origStyles["D51"] = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"

# --- Apply updated values row by row ---

$ws.Range("D2").Value = "20.576.62"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").Value = "1.475.47"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").Value = "0.9758"
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").Value = "278.41"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").Value = "0.3660"
$ws.Range("E7").Value = "  -1.56%  "

$ws.Range("D8").Value = "0.3067"
$ws.Range("E8").Value = "  -3.62%  "

$ws.Range("D9").Value = "40.24"
$ws.Range("E9").Value = "  -1.90%  "

$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("D11").Value = "0.06633"
$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "0.9986"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").Value = "5.484"
$ws.Range("E13").Value = "  -2.12%  "

$ws.Range("D14").Value = "18.12"
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").Value = "6.188"
$ws.Range("E15").Value = "  -1.53%  "

$ws.Range("D16").Value = "0.00001033"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").Value = "0.9754"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "1.477.52"
$ws.Range("E18").Value = "  +2.02%  "

$ws.Range("D19").Value = "0.05906"
$ws.Range("E19").Value = "  +1.91%  "

$ws.Range("D20").Value = "69.96"
$ws.Range("E20").Value = "  -3.33%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.474"
$ws.Range("E21").Value = "  -3.79%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "14.56"
$ws.Range("E22").Value = "  -2.62%  "

$ws.Range("D23").Value = "11.06"
$ws.Range("E23").Value = "  -1.51%  "

$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").Value = "20.595.67"
$ws.Range("E24").Value = "  +0.94%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.223"
$ws.Range("E25").Value = "  -2.88%  "

$ws.Range("D26").Value = "141.52"
$ws.Range("E26").Value = "  +3.99%  "

$ws.Range("D27").Value = "2.123"
$ws.Range("E27").Value = "  -9.53%  "

$ws.Range("D28").Value = "17.26"
$ws.Range("E28").Value = "  -1.81%  "

$ws.Range("D29").Value = "1.635.15"
$ws.Range("E29").Value = "  +1.77%  "

$ws.Range("D30").Value = "114.40"
$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("D31").Value = "3.885"
$ws.Range("E31").Value = "  -2.98%  "

$ws.Range("D32").Value = "0.8178"
$ws.Range("E32").Value = "  -4.10%  "

$ws.Range("D33").Value = "4.965"
$ws.Range("E33").Value = "  -7.72%  "

$ws.Range("D34").Value = "0.08018"
$ws.Range("E34").Value = "  +2.07%  "

$ws.Range("E35").Value = "  +2.40%  "

$ws.Range("D36").Value = "1.234"
$ws.Range("E36").Value = "  +9.81%  "

$ws.Range("D37").Value = "0.05807"
$ws.Range("E37").Value = "  -2.30%  "

$ws.Range("D38").Value = "4.729"
$ws.Range("E38").Value = "  -4.69%  "

$ws.Range("D39").Value = "0.9743"
$ws.Range("E39").Value = "  -1.63%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "7.685"
$ws.Range("E40").Value = "  -0.72%  "

$ws.Range("E41").Value = "  -2.53%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.02041"
$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("D43").Value = "0.1889"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("D44").Value = "0.5299"
$ws.Range("E44").Value = "  -2.63%  "

$ws.Range("D45").Value = "3.516"
$ws.Range("E45").Value = "  -1.67%  "

$ws.Range("D46").Value = "12.14"
$ws.Range("E46").Value = "  -3.57%  "

$ws.Range("D47").Value = "119.27"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("D48").Value = "0.5203"
$ws.Range("E48").Value = "  -3.45%  "

$ws.Range("D49").Value = "1.788"
$ws.Range("E49").Value = "  -3.10%  "

$ws.Range("D50").Value = "0.06462"
$ws.Range("E50").Value = "  +1.33%  "

$ws.Range("D51").Value = "0.9949"
$ws.Range("E51").Value = "  -0.13%  "

# Restore original styles on the price cells now that the text values are
# safely stored.
$ws.Range("D4").Style = $origStyles["D4"]
$ws.Range("D5").Style = $origStyles["D5"]
$ws.Range("D6").Style = $origStyles["D6"]
$ws.Range("D7").Style = $origStyles["D7"]
$ws.Range("D8").Style = $origStyles["D8"]
$ws.Range("D9").Style = $origStyles["D9"]
$ws.Range("D11").Style = $origStyles["D11"]
$ws.Range("D12").Style = $origStyles["D12"]
$ws.Range("D13").Style = $origStyles["D13"]
$ws.Range("D14").Style = $origStyles["D14"]
$ws.Range("D15").Style = $origStyles["D15"]
$ws.Range("D16").Style = $origStyles["D16"]
$ws.Range("D17").Style = $origStyles["D17"]
$ws.Range("D19").Style = $origStyles["D19"]
$ws.Range("D20").Style = $origStyles["D20"]
$ws.Range("D21").Style = $origStyles["D21"]
$ws.Range("D22").Style = $origStyles["D22"]
$ws.Range("D23").Style = $origStyles["D23"]
$ws.Range("D25").Style = $origStyles["D25"]
$ws.Range("D26").Style = $origStyles["D26"]
$ws.Range("D27").Style = $origStyles["D27"]
$ws.Range("D28").Style = $origStyles["D28"]
$ws.Range("D30").Style = $origStyles["D30"]
$ws.Range("D31").Style = $origStyles["D31"]
$ws.Range("D32").Style = $origStyles["D32"]
$ws.Range("D33").Style = $origStyles["D33"]
$ws.Range("D34").Style = $origStyles["D34"]
$ws.Range("D36").Style = $origStyles["D36"]
$ws.Range("D37").Style = $origStyles["D37"]
$ws.Range("D38").Style = $origStyles["D38"]
$ws.Range("D39").Style = $origStyles["D39"]
$ws.Range("D40").Style = $origStyles["D40"]
$ws.Range("D42").Style = $origStyles["D42"]
$ws.Range("D43").Style = $origStyles["D43"]
$ws.Range("D44").Style = $origStyles["D44"]
$ws.Range("D45").Style = $origStyles["D45"]
$ws.Range("D46").Style = $origStyles["D46"]
$ws.Range("D47").Style = $origStyles["D47"]
$ws.Range("D48").Style = $origStyles["D48"]
$ws.Range("D49").Style = $origStyles["D49"]
$ws.Range("D50").Style = $origStyles["D50"]
$ws.Range("D51").Style = $origStyles["D51"]
